# Update vm_pu.xlsx results for case with 380 kV done
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.037174652730593
$ws.Cells.Item(2, 4).Value = 1.039345150821084
$ws.Cells.Item(2, 5).Value = 1.044846267074777
$ws.Cells.Item(2, 6).Value = 1.052998016541885
$ws.Cells.Item(2, 9).Value = 1.039510812905085
$ws.Cells.Item(2, 10).Value = 1.042278957846756
$ws.Cells.Item(2, 11).Value = 1.042130602937207
$ws.Cells.Item(2, 12).Value = 1.047616169536049
$ws.Cells.Item(2, 13).Value = 1.055745200838477
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.038054295861997
$ws.Cells.Item(3, 4).Value = 1.039992116374041
$ws.Cells.Item(3, 5).Value = 1.045662683957637
$ws.Cells.Item(3, 6).Value = 1.053975934325667
$ws.Cells.Item(3, 9).Value = 1.039731913920884
$ws.Cells.Item(3, 10).Value = 1.042803222087647
$ws.Cells.Item(3, 11).Value = 1.042588198571586
$ws.Cells.Item(3, 12).Value = 1.048243893789712
$ws.Cells.Item(3, 13).Value = 1.056535652856418
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.038623885310594
$ws.Cells.Item(4, 4).Value = 1.040410971569566
$ws.Cells.Item(4, 5).Value = 1.046191759196989
$ws.Cells.Item(4, 6).Value = 1.054609856701185
$ws.Cells.Item(4, 9).Value = 1.039873779529956
$ws.Cells.Item(4, 10).Value = 1.043142209178983
$ws.Cells.Item(4, 11).Value = 1.042883832314429
$ws.Cells.Item(4, 12).Value = 1.048650220296053
$ws.Cells.Item(4, 13).Value = 1.057047652390872
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.038863435537513
$ws.Cells.Item(5, 4).Value = 1.04058711067115
$ws.Cells.Item(5, 5).Value = 1.046414372159535
$ws.Cells.Item(5, 6).Value = 1.054876630050585
$ws.Cells.Item(5, 9).Value = 1.039933131484263
$ws.Cells.Item(5, 10).Value = 1.043284659140218
$ws.Cells.Item(5, 11).Value = 1.043008005224954
$ws.Cells.Item(5, 12).Value = 1.048821073916485
$ws.Cells.Item(5, 13).Value = 1.057263021126956
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.038903662597333
$ws.Cells.Item(6, 4).Value = 1.040616688227943
$ws.Cells.Item(6, 5).Value = 1.046451760936925
$ws.Cells.Item(6, 6).Value = 1.054921438397865
$ws.Cells.Item(6, 9).Value = 1.039943080009421
$ws.Cells.Item(6, 10).Value = 1.043308573552787
$ws.Cells.Item(6, 11).Value = 1.043028847831878
$ws.Cells.Item(6, 12).Value = 1.048849762946944
$ws.Cells.Item(6, 13).Value = 1.057299189724819
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.038627085821486
$ws.Cells.Item(7, 4).Value = 1.040413324944993
$ws.Cells.Item(7, 5).Value = 1.046194733018695
$ws.Cells.Item(7, 6).Value = 1.05461342027138
$ws.Cells.Item(7, 9).Value = 1.039874573727554
$ws.Cells.Item(7, 10).Value = 1.043144112840303
$ws.Cells.Item(7, 11).Value = 1.042885491958488
$ws.Cells.Item(7, 12).Value = 1.048652503118592
$ws.Cells.Item(7, 13).Value = 1.057050529672424
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.037471848323503
$ws.Cells.Item(8, 4).Value = 1.03956374829936
$ws.Cells.Item(8, 5).Value = 1.04512201263919
$ws.Cells.Item(8, 6).Value = 1.053328270601945
$ws.Cells.Item(8, 9).Value = 1.039585783393737
$ws.Cells.Item(8, 10).Value = 1.042456185870477
$ws.Cells.Item(8, 11).Value = 1.042285344366882
$ws.Cells.Item(8, 12).Value = 1.047828280727888
$ws.Cells.Item(8, 13).Value = 1.056012228432013
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.035439304997625
$ws.Cells.Item(9, 4).Value = 1.038068478633945
$ws.Cells.Item(9, 5).Value = 1.043237928755869
$ws.Cells.Item(9, 6).Value = 1.05107250361547
$ws.Cells.Item(9, 9).Value = 1.039067726477223
$ws.Cells.Item(9, 10).Value = 1.041242127897351
$ws.Cells.Item(9, 11).Value = 1.041224323500111
$ws.Cells.Item(9, 12).Value = 1.046377074858833
$ws.Cells.Item(9, 13).Value = 1.054186684373468
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.034086459169528
$ws.Cells.Item(10, 4).Value = 1.037072929043329
$ws.Cells.Item(10, 5).Value = 1.041986123100413
$ws.Cells.Item(10, 6).Value = 1.049574687293532
$ws.Cells.Item(10, 9).Value = 1.038716230146761
$ws.Cells.Item(10, 10).Value = 1.040431581804161
$ws.Cells.Item(10, 11).Value = 1.040514698543595
$ws.Cells.Item(10, 12).Value = 1.045410469421601
$ws.Cells.Item(10, 13).Value = 1.052972476047896
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.033501195698953
$ws.Cells.Item(11, 4).Value = 1.036642171859087
$ws.Cells.Item(11, 5).Value = 1.041445104562645
$ws.Cells.Item(11, 6).Value = 1.048927565109512
$ws.Cells.Item(11, 9).Value = 1.038562584177112
$ws.Cells.Item(11, 10).Value = 1.040080341447528
$ws.Cells.Item(11, 11).Value = 1.04020689640788
$ws.Cells.Item(11, 12).Value = 1.04499213889361
$ws.Cells.Item(11, 13).Value = 1.052447396459219
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.03328388331554
$ws.Cells.Item(12, 4).Value = 1.036482219428597
$ws.Cells.Item(12, 5).Value = 1.041244301190728
$ws.Cells.Item(12, 6).Value = 1.04868741329845
$ws.Cells.Item(12, 9).Value = 1.038505296681528
$ws.Cells.Item(12, 10).Value = 1.039949835861339
$ws.Cells.Item(12, 11).Value = 1.040092486456272
$ws.Cells.Item(12, 12).Value = 1.044836785959507
$ws.Cells.Item(12, 13).Value = 1.052252462119796
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.033330493895324
$ws.Cells.Item(13, 4).Value = 1.036516527482143
$ws.Cells.Item(13, 5).Value = 1.041287367157046
$ws.Cells.Item(13, 6).Value = 1.048738916779682
$ws.Cells.Item(13, 9).Value = 1.038517594832068
$ws.Cells.Item(13, 10).Value = 1.039977831512886
$ws.Cells.Item(13, 11).Value = 1.040117031315896
$ws.Cells.Item(13, 12).Value = 1.044870108149251
$ws.Cells.Item(13, 13).Value = 1.052294271487684
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.033483230939567
$ws.Cells.Item(14, 4).Value = 1.036628949113092
$ws.Cells.Item(14, 5).Value = 1.041428502922487
$ws.Cells.Item(14, 6).Value = 1.048907709626312
$ws.Cells.Item(14, 9).Value = 1.038557853188385
$ws.Cells.Item(14, 10).Value = 1.04006955461016
$ws.Cells.Item(14, 11).Value = 1.040197440841462
$ws.Cells.Item(14, 12).Value = 1.044979296678231
$ws.Cells.Item(14, 13).Value = 1.052431281004755
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.033577347981655
$ws.Cells.Item(15, 4).Value = 1.03669822245736
$ws.Cells.Item(15, 5).Value = 1.041515481910429
$ws.Cells.Item(15, 6).Value = 1.049011737426378
$ws.Cells.Item(15, 9).Value = 1.038582629016477
$ws.Cells.Item(15, 10).Value = 1.040126063064768
$ws.Cells.Item(15, 11).Value = 1.040246973429563
$ws.Cells.Item(15, 12).Value = 1.045046575829429
$ws.Cells.Item(15, 13).Value = 1.052515710858158
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.034125312402989
$ws.Cells.Item(16, 4).Value = 1.037101523901911
$ws.Cells.Item(16, 5).Value = 1.042022050372243
$ws.Cells.Item(16, 6).Value = 1.049617665175446
$ws.Cells.Item(16, 9).Value = 1.038726396720381
$ws.Cells.Item(16, 10).Value = 1.040454886902949
$ws.Cells.Item(16, 11).Value = 1.040535115285694
$ws.Cells.Item(16, 12).Value = 1.045438237294213
$ws.Cells.Item(16, 13).Value = 1.053007338305484
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.034469178167397
$ws.Cells.Item(17, 4).Value = 1.037354591632897
$ws.Cells.Item(17, 5).Value = 1.042340081963366
$ws.Cells.Item(17, 6).Value = 1.049998134819792
$ws.Cells.Item(17, 9).Value = 1.038816191806088
$ws.Cells.Item(17, 10).Value = 1.040661078264535
$ws.Cells.Item(17, 11).Value = 1.040715717959829
$ws.Cells.Item(17, 12).Value = 1.045683974859594
$ws.Cells.Item(17, 13).Value = 1.053315906385237
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.034669800081639
$ws.Cells.Item(18, 4).Value = 1.037502232758793
$ws.Cells.Item(18, 5).Value = 1.042525682874575
$ws.Cells.Item(18, 6).Value = 1.050220195373043
$ws.Cells.Item(18, 9).Value = 1.038868428222873
$ws.Cells.Item(18, 10).Value = 1.040781320214502
$ws.Cells.Item(18, 11).Value = 1.040821009249829
$ws.Cells.Item(18, 12).Value = 1.045827330184911
$ws.Cells.Item(18, 13).Value = 1.053495954506313
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034738215553881
$ws.Cells.Item(19, 4).Value = 1.037552579778685
$ws.Cells.Item(19, 5).Value = 1.04258898464396
$ws.Cells.Item(19, 6).Value = 1.05029593582192
$ws.Cells.Item(19, 9).Value = 1.038886215802908
$ws.Cells.Item(19, 10).Value = 1.040822315172505
$ws.Cells.Item(19, 11).Value = 1.040856902142757
$ws.Cells.Item(19, 12).Value = 1.0458762141312
$ws.Cells.Item(19, 13).Value = 1.053557357345557
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.034432279358409
$ws.Cells.Item(20, 4).Value = 1.037327436642795
$ws.Cells.Item(20, 5).Value = 1.042305949987785
$ws.Cells.Item(20, 6).Value = 1.049957299634496
$ws.Cells.Item(20, 9).Value = 1.038806572071416
$ws.Cells.Item(20, 10).Value = 1.04063895856261
$ws.Cells.Item(20, 11).Value = 1.040696346289629
$ws.Cells.Item(20, 12).Value = 1.045657607388401
$ws.Cells.Item(20, 13).Value = 1.053282793157513
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.033438251453894
$ws.Cells.Item(21, 4).Value = 1.036595842353146
$ws.Cells.Item(21, 5).Value = 1.041386937666245
$ws.Cells.Item(21, 6).Value = 1.048857998305008
$ws.Cells.Item(21, 9).Value = 1.03854600407812
$ws.Cells.Item(21, 10).Value = 1.040042545519859
$ws.Cells.Item(21, 11).Value = 1.040173764399655
$ws.Cells.Item(21, 12).Value = 1.044947142439367
$ws.Cells.Item(21, 13).Value = 1.052390932246224
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.032813733190111
$ws.Cells.Item(22, 4).Value = 1.036136149840035
$ws.Cells.Item(22, 5).Value = 1.040810015690758
$ws.Cells.Item(22, 6).Value = 1.048168087007255
$ws.Cells.Item(22, 9).Value = 1.038380922127256
$ws.Cells.Item(22, 10).Value = 1.039667330242526
$ws.Cells.Item(22, 11).Value = 1.039844742583083
$ws.Cells.Item(22, 12).Value = 1.044500640210578
$ws.Cells.Item(22, 13).Value = 1.051830783950876
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.033144757498676
$ws.Cells.Item(23, 4).Value = 1.036379813551688
$ws.Cells.Item(23, 5).Value = 1.041115767360328
$ws.Cells.Item(23, 6).Value = 1.048533701808991
$ws.Cells.Item(23, 9).Value = 1.038468553694208
$ws.Cells.Item(23, 10).Value = 1.039866260094642
$ws.Cells.Item(23, 11).Value = 1.040019205935006
$ws.Cells.Item(23, 12).Value = 1.044737320579427
$ws.Cells.Item(23, 13).Value = 1.052127671833456
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.034448952184954
$ws.Cells.Item(24, 4).Value = 1.037339706718344
$ws.Cells.Item(24, 5).Value = 1.042321372454288
$ws.Cells.Item(24, 6).Value = 1.049975750867528
$ws.Cells.Item(24, 9).Value = 1.03881091924665
$ws.Cells.Item(24, 10).Value = 1.040648953584171
$ws.Cells.Item(24, 11).Value = 1.040705099671844
$ws.Cells.Item(24, 12).Value = 1.045669521649638
$ws.Cells.Item(24, 13).Value = 1.053297755397285
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035964387020788
$ws.Cells.Item(25, 4).Value = 1.038454819786582
$ws.Cells.Item(25, 5).Value = 1.043724267111975
$ws.Cells.Item(25, 6).Value = 1.051654617580935
$ws.Cells.Item(25, 9).Value = 1.039202739219382
$ws.Cells.Item(25, 10).Value = 1.041556202679
$ws.Cells.Item(25, 11).Value = 1.041499029033673
$ws.Cells.Item(25, 12).Value = 1.046752099165962
$ws.Cells.Item(25, 13).Value = 1.054658139976881
